$d = $word.ActiveDocument

# --- The "Meta description" paragraph currently sits right after the
# Heading1 title (paragraph 2). We need to move its content (an empty run
# + a bold "Meta description" run + a plain ": Want to play..." run) down
# to just before the final (italic) paragraph, turning it into a new bold
# title paragraph, and remove it from its original spot. ---
$metaPara = $d.Paragraphs.Item(2)

# Insert a brand new empty paragraph right before the final paragraph
# FIRST (while the meta-description paragraph is still intact), so that
# nothing has shifted yet when we go fetch its formatted content.
$priorPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$priorPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)

# Make sure the freshly inserted paragraph doesn't inherit list/heading
# formatting from its neighbour.
$newPara.Style = "Normal"
$newPara.Range.ParagraphFormat.Reset() | Out-Null

# Grab the meta-description paragraph's formatted content (fresh, right
# before pasting) and copy it into the new paragraph.
$metaRange = $d.Range($metaPara.Range.Start, $metaPara.Range.End - 1)
$destRange = $d.Range($newPara.Range.Start, $newPara.Range.Start)
$destRange.FormattedText = $metaRange.FormattedText

# Now remove the original "Meta description: ..." paragraph entirely.
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# --- Turn the copied paragraph into the new bold title paragraph. ---
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)

# Drop the ": Want to play ..." trailing run entirely.
$find = $newPara.Range.Duplicate
$find.Find.Execute(": Want to play Bear Money slot for free? Read our review before you start and find out all the pros and cons of the game. Discover the Multi Cash Collector bonus and attractive symbols. ", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$find.Delete()

# Turn the remaining "Meta description" bold run into the new title text.
$find2 = $newPara.Range.Duplicate
$find2.Find.Execute("Meta description", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Play Bear Money Slot for Free - Review & Rating 2021", 2) | Out-Null

# --- Replace the text of the final (italic) paragraph with the old meta
# description wording, keeping its italic formatting intact. ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$findLast = $lastPara.Range.Duplicate
$oldImagePrompt = 'Create a cartoon style feature image for "Bear Money" that features a happy Maya warrior with glasses. The image should showcase the forest setting of the game, with trees and wildlife in the background. The Maya warrior should be holding a picnic basket and a jar of honey, with a big smile on their face as they outsmart the thieving bears. The colors should be bright and eye-catching, with a playful and fun tone that captures the spirit of the game.'
$newMetaText = 'Want to play Bear Money slot for free? Read our review before you start and find out all the pros and cons of the game. Discover the Multi Cash Collector bonus and attractive symbols. '
$findLast.Find.Execute($oldImagePrompt, $true, $false, $false, $false, $false, `
    $true, 1, $false, $newMetaText, 2) | Out-Null
